# Commit: "Used PR#8 and added code for Home & DSIntro pages."
# Adds six new worksheets (Sheet5..Sheet10) holding DSA-topic option/code/result
# tables, mirroring the PR that introduced the LinkedList & Stack "Home" pages.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the six new sheets, appended after Sheet4 (so they land at the
#    end of the tab strip as Sheet5..Sheet10 with sheetId 6..11).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet5  = $wb.Worksheets.Add($null, $lastSheet)
$sheet6  = $wb.Worksheets.Add($null, $sheet5)
$sheet7  = $wb.Worksheets.Add($null, $sheet6)
$sheet8  = $wb.Worksheets.Add($null, $sheet7)
$sheet9  = $wb.Worksheets.Add($null, $sheet8)
$sheet10 = $wb.Worksheets.Add($null, $sheet9)

# Sheet5 / Sheet6 stay blank (placeholders), just leave default selection.
$sheet5.Range("A1").Select() | Out-Null
$sheet6.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Sheet7 - "OptionOnLinkedList" / code / errorMessage, undefined-name demo
# ---------------------------------------------------------------------------
$sheet7.Cells.Item(1,1).Value = "OptionOnLinkedList"
$sheet7.Cells.Item(1,2).Value = "code"
$sheet7.Cells.Item(1,3).Value = "errorMessage"

$sheet7.Cells.Item(2,1).Value = "Introduction"
$sheet7.Cells.Item(3,1).Value = "Creating Linked LIst"
$sheet7.Cells.Item(4,1).Value = "Types of Linked List"
$sheet7.Cells.Item(5,1).Value = "Implement Linked List in Python"
$sheet7.Cells.Item(6,1).Value = "Traversal"
$sheet7.Cells.Item(7,1).Value = "Insertion"
$sheet7.Cells.Item(8,1).Value = "Deletion"

for ($r = 2; $r -le 8; $r++) {
    $sheet7.Cells.Item($r,2).Value = "abc123"
    $sheet7.Cells.Item($r,3).Value = "NameError: name 'abc123' is not defined on line 1"
}

$sheet7.Columns.Item(1).ColumnWidth = 28.333333333333336
$sheet7.Columns.Item(3).ColumnWidth = 44.666666666666664
$sheet7.Range("C2").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Sheet8 - same option list, working "print(...)" snippets this time
# ---------------------------------------------------------------------------
$sheet8.Cells.Item(1,1).Value = "OptionOnLinkedList"
$sheet8.Cells.Item(1,2).Value = "code"
$sheet8.Cells.Item(1,3).Value = "errorMessage"

$sheet8.Cells.Item(2,1).Value = "Introduction"
$sheet8.Cells.Item(3,1).Value = "Creating Linked LIst"
$sheet8.Cells.Item(4,1).Value = "Types of Linked List"
$sheet8.Cells.Item(5,1).Value = "Implement Linked List in Python"
$sheet8.Cells.Item(6,1).Value = "Traversal"
$sheet8.Cells.Item(7,1).Value = "Insertion"
$sheet8.Cells.Item(8,1).Value = "Deletion"

$sheet8.Cells.Item(2,2).Value = "print('Hello World! Introduction')"
$sheet8.Cells.Item(3,2).Value = "print('Hello World! Creating Linked LIst')"
$sheet8.Cells.Item(4,2).Value = "print('Hello World! Types of Linked List')"
$sheet8.Cells.Item(5,2).Value = "print('Hello World! Implement Linked List in Python')"
$sheet8.Cells.Item(6,2).Value = "print('Hello World! Traversal')"
$sheet8.Cells.Item(7,2).Value = "print('Hello World! Insertion')"
$sheet8.Cells.Item(8,2).Value = "print('Hello World! Deletion')"

$sheet8.Cells.Item(2,3).Value = "Hello World! Introduction"
$sheet8.Cells.Item(3,3).Value = "Hello World! Creating Linked LIst"
$sheet8.Cells.Item(4,3).Value = "Hello World! Types of Linked List"
$sheet8.Cells.Item(5,3).Value = "Hello World! Implement Linked List in Python"
$sheet8.Cells.Item(6,3).Value = "Hello World! Traversal"
$sheet8.Cells.Item(7,3).Value = "Hello World! Insertion"
$sheet8.Cells.Item(8,3).Value = "Hello World! Deletion"

$sheet8.Columns.Item(1).ColumnWidth = 28.333333333333336
$sheet8.Columns.Item(2).ColumnWidth = 46.83333333333333
$sheet8.Columns.Item(3).ColumnWidth = 40.166666666666664
$sheet8.Range("B4").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Sheet9 - "OptionOnStack" / code / errorMessage, undefined-name demo
# ---------------------------------------------------------------------------
$sheet9.Cells.Item(1,1).Value = "OptionOnStack"
$sheet9.Cells.Item(1,2).Value = "code"
$sheet9.Cells.Item(1,3).Value = "errorMessage"

$sheet9.Cells.Item(2,1).Value = "Operations in Stack"
$sheet9.Cells.Item(3,1).Value = "Implementation"
$sheet9.Cells.Item(4,1).Value = "Applications"

for ($r = 2; $r -le 4; $r++) {
    $sheet9.Cells.Item($r,2).Value = "Abcd"
    $sheet9.Cells.Item($r,3).Value = "NameError: name 'Abcd' is not defined on line 1"
}

$sheet9.Columns.Item(1).ColumnWidth = 28.333333333333336
$sheet9.Columns.Item(3).ColumnWidth = 44.666666666666664
$sheet9.Range("A4").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Sheet10 - "OptionOnStack" / code / message, working "print('Hello World')"
# ---------------------------------------------------------------------------
$sheet10.Cells.Item(1,1).Value = "OptionOnStack"
$sheet10.Cells.Item(1,2).Value = "code"
$sheet10.Cells.Item(1,3).Value = "message"

$sheet10.Cells.Item(2,1).Value = "Operations in Stack"
$sheet10.Cells.Item(3,1).Value = "Implementation"
$sheet10.Cells.Item(4,1).Value = "Applications"

for ($r = 2; $r -le 4; $r++) {
    $sheet10.Cells.Item($r,2).Value = "print('Hello World')"
    $sheet10.Cells.Item($r,3).Value = "Hello World"
}

$sheet10.Columns.Item(1).ColumnWidth = 17.666666666666668
$sheet10.Columns.Item(2).ColumnWidth = 17.166666666666668
$sheet10.Columns.Item(3).ColumnWidth = 10.5
$sheet10.Range("C5").Select() | Out-Null

# Sheet10 (the last-created sheet) stays the active tab, matching the
# workbook's new activeTab/tabSelected state.
$sheet10.Activate() | Out-Null
